$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.975.59"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "1.648.48"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.78%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.25"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06431"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.326"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.94%  "

$ws.Range("D13").Value = "1.652.68"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5474"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "0.0₅7908"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("E16").Value = "  +1.81%  "

$ws.Range("D17").Value = "26.043.21"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.488"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.872"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1151"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.910"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.98%  "

$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.243"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05025"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.285"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.367"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8956"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.598"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.60%  "

$ws.Range("D37").Value = "1.136.33"
$ws.Range("E37").Value = "  -3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5543"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01566"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.007"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.551"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.660"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8184"
$ws.Range("D43").Style = "Normal"

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈125"
$ws.Range("E45").Value = "  +9.30%  "

$ws.Range("D46").Value = "1.785.93"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4546"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05096"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "

